# Insert a new record row at row 462, shifting the existing rows 462-558
# down to 463-559, and populate the new row with the new data record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 462 (pushes 462..558 down to 463..559)
$ws.Rows.Item(462).Insert()

# Populate the newly inserted row 462 with the new data record.
$ws.Range("A462").Value2 = 7
$ws.Range("B462").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C462").Value2 = "Ñuble"
$ws.Range("D462").Value2 = 45258
$ws.Range("E462").Value2 = 16
$ws.Range("F462").Value2 = 100114013
$ws.Range("G462").Value2 = "Zanahoria"
$ws.Range("H462").Value2 = "Sin especificar"
$ws.Range("I462").Value2 = "Primera"
$ws.Range("J462").Value2 = 150
$ws.Range("K462").Value2 = 7000
$ws.Range("L462").Value2 = 7000
$ws.Range("M462").Value2 = 7000
$ws.Range("N462").Value2 = "$/saco 20 kilos"
$ws.Range("O462").Value2 = "Provincia de Diguillín"
$ws.Range("P462").Value2 = 350
$ws.Range("Q462").Value2 = 20
$ws.Range("R462").Value2 = "Hortaliza"
